$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17 (Svi Mykhailiuk) is being folded into row 16 (which currently holds
# Reggie Jackson, who is being dropped from the roster). Copy row 17's player
# data into row 16 -- but leave row 16's existing hyperlink relationship
# (K16) untouched, only its display text changes -- then delete row 17.

$ws.Range("B16").Value = 10
$ws.Range("C16").Value = "Svi Mykhailiuk"
$ws.Range("D16").Value = "SF"
$ws.Range("E16").Value = "6-7"
$ws.Range("F16").Value = 205
$ws.Range("G16").Value = "June 10, 1997"
$ws.Range("H16").Value = "ua"
$ws.Range("I16").Value = "4"
$ws.Range("J16").Value = "Kansas"
$ws.Range("K16").Value = "https://www.basketball-reference.com/players/m/mykhasv01.html"

$ws.Rows("17").Delete()
